$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 981.7
$ws.Range("I19").Value = 964
$ws.Range("J19").Value = 986.125
$ws.Range("K19").Value = 964
$ws.Range("L19").Value = 986.125
$ws.Range("M19").Value = -789
$ws.Range("N19").Value = -1336.125

$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("L64").ClearContents()
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = 0

$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("L67").ClearContents()
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = 0

$ws.Range("H100").Value = 2437.5
$ws.Range("I100").Value = 2437.5
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2437.5
$ws.Range("L100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -1896.5

$ws.Range("H113").Value = 2962.2
$ws.Range("I113").Value = 2269.6667
$ws.Range("J113").Value = 4001
$ws.Range("K113").Value = 2269.6667
$ws.Range("L113").Value = 4001
$ws.Range("M113").Value = 984.3332999999998
$ws.Range("N113").Value = -10509

$ws.Range("H116").Value = 3712.2856
$ws.Range("I116").Value = 1990
$ws.Range("J116").Value = 4401.2
$ws.Range("K116").Value = 1990
$ws.Range("L116").Value = 4401.2
$ws.Range("M116").Value = 1452
$ws.Range("N116").Value = -11285.2

$ws.Range("H129").Value = 756.2
$ws.Range("J129").Value = 811.55554
$ws.Range("L129").Value = 2434.66662
$ws.Range("N129").Value = -12434.66662

$ws.Range("H137").Value = 1492.7142
$ws.Range("I137").Value = 1500.3636
$ws.Range("J137").Value = 1487.7646
$ws.Range("K137").Value = 4501.0908
$ws.Range("L137").Value = 4463.293799999999
$ws.Range("M137").Value = -1951.0908
$ws.Range("N137").Value = -9563.293799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1063.091
$ws.Range("I2").Value = 964.8
$ws.Range("K2").Value = 964.8
$ws.Range("M2").Value = -851.8

$ws.Range("H32").Value = 5143.6
$ws.Range("I32").Value = 5143.6
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5143.6
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -4856.6

$ws.Range("H61").Value = 250002130
$ws.Range("I61").Value = 250002130
$ws.Range("K61").Value = 250002130
$ws.Range("M61").Value = -250001918

$ws.Range("H63").Value = 2166.5881
$ws.Range("I63").Value = 2072.9092
$ws.Range("J63").Value = 2338.3333
$ws.Range("K63").Value = 2072.9092
$ws.Range("L63").Value = 2338.3333
$ws.Range("M63").Value = -1386.9092
$ws.Range("N63").Value = -3710.3333

$ws.Range("H66").Value = 2166.5881
$ws.Range("I66").Value = 2072.9092
$ws.Range("J66").Value = 2338.3333
$ws.Range("K66").Value = 10364.546
$ws.Range("L66").Value = 11691.6665
$ws.Range("M66").Value = -6932.546
$ws.Range("N66").Value = -18555.6665

$ws.Range("H116").Value = 1063.091
$ws.Range("I116").Value = 964.8
$ws.Range("K116").Value = 964.8
$ws.Range("M116").Value = 1329.2

$ws.Range("H132").Value = 3093.8235
$ws.Range("I132").Value = 2609.0908
$ws.Range("K132").Value = 7827.2724
$ws.Range("M132").Value = -5297.2724

$ws.Range("H136").Value = 250002130
$ws.Range("I136").Value = 250002130
$ws.Range("K136").Value = 750006390
$ws.Range("M136").Value = -750003840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1063.091
$ws.Range("I3").Value = 964.8
$ws.Range("K3").Value = 964.8
$ws.Range("M3").Value = -850.8

$ws.Range("H94").Value = 25000846
$ws.Range("I94").Value = 31250634
$ws.Range("K94").Value = 31250634
$ws.Range("M94").Value = -31250183

$ws.Range("H105").Value = 112210904
$ws.Range("I105").Value = 126237110
$ws.Range("K105").Value = 126237110
$ws.Range("M105").Value = -126235363

$ws.Range("H134").Value = 6348.524
$ws.Range("I134").Value = 1077.6471
$ws.Range("K134").Value = 3232.9413
$ws.Range("M134").Value = -697.9412999999995

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 66667976
$ws.Range("I16").Value = 90910360
$ws.Range("J16").Value = 1420
$ws.Range("K16").Value = 90910360
$ws.Range("L16").Value = 1420
$ws.Range("M16").Value = -90910073
$ws.Range("N16").Value = -1994

$ws.Range("H94").Value = 2179.6
$ws.Range("I94").Value = 2999
$ws.Range("J94").Value = 1974.75
$ws.Range("K94").Value = 2999
$ws.Range("L94").Value = 1974.75
$ws.Range("M94").Value = -2548
$ws.Range("N94").Value = -2876.75

$ws.Range("H113").Value = 66667976
$ws.Range("I113").Value = 90910360
$ws.Range("J113").Value = 1420
$ws.Range("K113").Value = 90910360
$ws.Range("L113").Value = 1420
$ws.Range("M113").Value = -90908190
$ws.Range("N113").Value = -5760

$ws.Range("H132").Value = 9995.77
$ws.Range("I132").Value = 21544.8
$ws.Range("J132").Value = 2777.625
$ws.Range("K132").Value = 64634.39999999999
$ws.Range("L132").Value = 8332.875
$ws.Range("M132").Value = -62104.39999999999
$ws.Range("N132").Value = -13392.875

$ws.Range("H134").Value = 33336680
$ws.Range("I134").Value = 4355.778
$ws.Range("J134").Value = 83335170
$ws.Range("K134").Value = 13067.334
$ws.Range("L134").Value = 250005510
$ws.Range("M134").Value = -10532.334
$ws.Range("N134").Value = -250010580

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 2757.5
$ws.Range("I126").Value = 2757.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 8272.5
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -3332.5

$ws.Range("H131").Value = 15875901
$ws.Range("I131").Value = 125000310
$ws.Range("J131").Value = 3259.5637
$ws.Range("K131").Value = 375000930
$ws.Range("L131").Value = 9778.6911
$ws.Range("M131").Value = -374995890
$ws.Range("N131").Value = -19858.6911

$ws.Range("H137").Value = 30007816
$ws.Range("I137").Value = 150003200
$ws.Range("J137").Value = 8971.200000000001
$ws.Range("K137").Value = 450009600
$ws.Range("L137").Value = 26913.6
$ws.Range("M137").Value = -450004500
$ws.Range("N137").Value = -37113.60000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1753.3334
$ws.Range("I102").Value = 1101.5
$ws.Range("K102").Value = 1101.5
$ws.Range("M102").Value = 520.5

$ws.Range("H132").Value = 3546.2354
$ws.Range("I132").Value = 3377.7144
$ws.Range("K132").Value = 10133.1432
$ws.Range("M132").Value = -7603.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2563.3928
$ws.Range("I40").Value = 1761.5
$ws.Range("K40").Value = 1761.5
$ws.Range("M40").Value = -1625.5

$ws.Range("H46").Value = 3614.2856
$ws.Range("I46").Value = 687.5
$ws.Range("J46").Value = 4785
$ws.Range("K46").Value = 687.5
$ws.Range("L46").Value = 4785
$ws.Range("M46").Value = -499.5
$ws.Range("N46").Value = -5161

$ws.Range("H61").Value = 1107.5
$ws.Range("I61").Value = 1041.6666
$ws.Range("K61").Value = 1041.6666
$ws.Range("M61").Value = -839.6666

$ws.Range("H68").Value = 1951.625
$ws.Range("J68").Value = 2999
$ws.Range("L68").Value = 2999
$ws.Range("N68").Value = -4497

$ws.Range("H71").Value = 1951.625
$ws.Range("J71").Value = 2999
$ws.Range("L71").Value = 14995
$ws.Range("N71").Value = -22483

$ws.Range("H82").Value = 2199.1538
$ws.Range("I82").Value = 2065.6667
$ws.Range("J82").Value = 2499.5
$ws.Range("K82").Value = 2065.6667
$ws.Range("L82").Value = 2499.5
$ws.Range("M82").Value = -1704.6667
$ws.Range("N82").Value = -3221.5

$ws.Range("H85").Value = 2199.1538
$ws.Range("I85").Value = 2065.6667
$ws.Range("J85").Value = 2499.5
$ws.Range("K85").Value = 2065.6667
$ws.Range("L85").Value = 2499.5
$ws.Range("M85").Value = -817.6667000000002
$ws.Range("N85").Value = -4995.5

$ws.Range("H93").Value = 1020.4
$ws.Range("I93").Value = 1001
$ws.Range("J93").Value = 1033.3334
$ws.Range("K93").Value = 1001
$ws.Range("L93").Value = 1033.3334
$ws.Range("M93").Value = 247
$ws.Range("N93").Value = -3529.3334

$ws.Range("H100").Value = 1601
$ws.Range("J100").Value = 2500
$ws.Range("L100").Value = 2500
$ws.Range("N100").Value = -3582

$ws.Range("H113").Value = 1107.5
$ws.Range("I113").Value = 1041.6666
$ws.Range("K113").Value = 1041.6666
$ws.Range("M113").Value = 1128.3334

$ws.Range("H122").Value = 50028800
$ws.Range("J122").Value = 21152.5
$ws.Range("L122").Value = 63457.5
$ws.Range("N122").Value = -68357.5

$ws.Range("H131").Value = 38000
$ws.Range("J131").Value = 38000
$ws.Range("L131").Value = 38000
$ws.Range("N131").Value = -48080

$ws.Range("H132").Value = 81422.14
$ws.Range("I132").Value = 28300.75
$ws.Range("J132").Value = 102670.7
$ws.Range("K132").Value = 84902.25
$ws.Range("L132").Value = 308012.1
$ws.Range("M132").Value = -82372.25
$ws.Range("N132").Value = -313072.1

$ws.Range("H136").Value = 15786.857
$ws.Range("I136").Value = 25752
$ws.Range("K136").Value = 77256
$ws.Range("M136").Value = -74706

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 495.33334
$ws.Range("I113").Value = 309.66666
$ws.Range("K113").Value = 928.9999799999999
$ws.Range("M113").Value = 1241.00002

$ws.Range("H122").Value = 10418573
$ws.Range("I122").Value = 10871485
$ws.Range("K122").Value = 32614455
$ws.Range("M122").Value = -32612005

$ws.Range("H129").Value = 75000
$ws.Range("J129").Value = 75000
$ws.Range("L129").Value = 75000
$ws.Range("N129").Value = -85000

$ws.Range("H132").Value = 6317
$ws.Range("I132").Value = 8468
$ws.Range("J132").Value = 4166
$ws.Range("K132").Value = 25404
$ws.Range("L132").Value = 12498
$ws.Range("M132").Value = -22874
$ws.Range("N132").Value = -17558

$ws.Range("H136").Value = 1144.7778
$ws.Range("I136").Value = 1141.8125
$ws.Range("J136").Value = 1149.091
$ws.Range("K136").Value = 3425.4375
$ws.Range("L136").Value = 3447.273
$ws.Range("M136").Value = -875.4375
$ws.Range("N136").Value = -8547.272999999999
